$d = $word.ActiveDocument

# The bibliography ends with "Janeiro: Editora Interciência , 2004." and is
# followed by a blank paragraph, then the site-chrome paragraphs
# "Ver no Jupiter Salvar em pdf Salvar em docx" and
# "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages.
#  Original theme under Creative Commons Attribution" which the site
# rebuild removed. Locate the "Janeiro..." paragraph, then delete the
# three paragraphs that directly follow it (the blank one plus the two
# site-chrome paragraphs), leaving the "Janeiro..." paragraph directly
# followed by the pre-existing blank paragraph that comes before the
# trailing page-break paragraph.

$targetIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Janeiro: Editora Interciência , 2004.") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -ne $null) {
    # Remove the paragraph right after "Janeiro..." three times in a row
    # (the collection re-indexes after each delete, so the paragraph that
    # needs removing is always back at $targetIndex + 1).
    for ($count = 1; $count -le 3; $count++) {
        $p = $d.Paragraphs.Item($targetIndex + 1)
        $r = $d.Range($p.Range.Start, $p.Range.End)
        $r.Delete()
    }
}
